# Rename the embedded logo pictures (Pearson logo in both footers, BTec logo
# in the header) so their shape "Name" matches the new export numbering:
#   - footer Pearson logos: image1.png -> image2.png
#   - header BTec logo:     image2.jpg -> image1.jpg
#
# Word exposes the picture's OOXML <wp:docPr>/name (and the paired
# <pic:cNvPr>/name) through InlineShape.Name, same as Shape.Name for a
# floating shape, so we just walk every section's headers/footers and rename
# each inline picture by matching on its (unique, unchanged) AlternativeText.

$d = $word.ActiveDocument

function Rename-LogoInHeaderFooter($headerFooter, $altText, $newName) {
    if ($headerFooter.Exists) {
        foreach ($ishp in $headerFooter.Range.InlineShapes) {
            if ($ishp.AlternativeText -eq $altText) {
                $ishp.Name = $newName
            }
        }
    }
}

foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        Rename-LogoInHeaderFooter $hf "BTec_Logo-Orange" "image1.jpg"
    }
    foreach ($hf in $sec.Footers) {
        Rename-LogoInHeaderFooter $hf "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" "image2.png"
    }
}
